$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exclude "中国美术学院" (row 2) from the admissions list. Its admission-
# guide cell (C2) carries a hyperlink, which Excel does not always retarget
# automatically on a row delete, so drop it explicitly first.
$ws.Range("C2").Hyperlinks.Delete()

# Delete the entire row so every following row shifts up by one.
$ws.Rows.Item(2).Delete()
